# Auto-generated edit script: update "remaining tickets" (F column) and
# "lowest price" (G column) values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 981
$ws.Range("F3").Value = 1046
$ws.Range("F5").Value = 907
$ws.Range("F6").Value = 493
$ws.Range("F7").Value = 493
$ws.Range("F8").Value = 770
$ws.Range("F9").Value = 181
$ws.Range("F10").Value = 1365
$ws.Range("F11").Value = 794
$ws.Range("F12").Value = 445
$ws.Range("F13").Value = 601
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 21
$ws.Range("F17").Value = 200
$ws.Range("F18").Value = 97
$ws.Range("F19").Value = 1480
$ws.Range("F20").Value = 172
$ws.Range("F21").Value = 23
$ws.Range("F22").Value = 449
$ws.Range("F23").Value = 44
$ws.Range("F24").Value = 385
$ws.Range("F26").Value = 622
$ws.Range("F27").Value = 10
$ws.Range("F28").Value = 190
$ws.Range("F29").Value = 696
$ws.Range("F31").Value = 1294
$ws.Range("F32").Value = 84
$ws.Range("F33").Value = 22
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 357
$ws.Range("G2").Value = "不可售"
$ws.Range("F3").Value = 118
$ws.Range("F6").Value = 193
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 680
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 357
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 981
$ws.Range("F5").Value = 1046
$ws.Range("F7").Value = 907
$ws.Range("F8").Value = 493
$ws.Range("F9").Value = 493
$ws.Range("F10").Value = 770
$ws.Range("F11").Value = 181
$ws.Range("F12").Value = 1365
$ws.Range("F13").Value = 794
$ws.Range("F14").Value = 118
$ws.Range("F16").Value = 445
$ws.Range("F17").Value = 601
$ws.Range("F19").Value = 148
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 21
$ws.Range("F22").Value = 200
$ws.Range("F23").Value = 97
$ws.Range("F24").Value = 1480
$ws.Range("F25").Value = 193
$ws.Range("F26").Value = 172
$ws.Range("F27").Value = 23
$ws.Range("F28").Value = 449
$ws.Range("F29").Value = 44
$ws.Range("F30").Value = 385
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 680
$ws.Range("F35").Value = 622
$ws.Range("F40").Value = 10
$ws.Range("F41").Value = 190
$ws.Range("F42").Value = 696
$ws.Range("F44").Value = 1294
$ws.Range("F45").Value = 84
$ws.Range("F46").Value = 22
